# Apply the "rel_migrace" edit:
#  - column I on sheet "src" switches from the uap2016!X ratio column to a
#    freshly computed E/C ratio (relative migration balance)
#  - the header label for column I changes from "zastavenost" to "rel_migrace"
#  - view/selection state on both "src" and "uap2016" sheets is updated to
#    reflect where the author left the cursor after the edit

$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("src")
$uap = $wb.Worksheets.Item("uap2016")

# Rename the column header (shared string table is rebuilt automatically on
# save, so this also takes care of removing the now-unused "zastavenost"
# string and appending the new "rel_migrace" string).
$src.Range("I1").Value2 = "rel_migrace"

# Recompute column I for every data row: relative migration balance =
# pristehovani (E) / pocet obyvatel (C) on sheet uap2016, row offset by 5.
for ($r = 2; $r -le 207; $r++) {
    $uapRow = $r + 5
    $src.Cells.Item($r, 9).Formula = "='uap2016'!E$uapRow/'uap2016'!C$uapRow"
}

# Update view/selection on the "uap2016" sheet: scroll the frozen pane back
# to the top (C7) and collapse the selection down to a single cell (E7).
$uap.Activate()
$uap.Range("E7").Select()

# Update view/selection on the "src" sheet: drop the scrolled-down
# topLeftCell and select the whole newly-computed column I instead, leaving
# "src" as the active sheet/tab.
$src.Activate()
$src.Range("I2:I207").Select()
